# Natmi following Dr Hou advice
#
# The sending/target-cluster cross table is rebuilt: every "Sending cluster"
# (ECs, FAPs, M1, M2, Neutro, sCs) now gets a row against each of the three
# possible "Target cluster" values (ECs, FAPs, sCs) -- 6 x 3 = 18 data rows
# (rows 2-19) instead of the original 12, and every numeric column is
# recomputed accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol,
#          D Target cluster, E..T the recomputed numeric metrics.
$data = @(
    @("ECs",    "Thbs1", "Tnfrsf11b", "ECs",  3, 1, 15.18897966666667,  45.566939,   0.01327029680642083, 0.01327029680642082, 1, 0.3333333333333333, 0.07580833333333332, 0.227425, 0.01036056733223033, 0.01036056733223033, 1.151451233563889,  10.363061102075,    0.0001374878035816041, 0.0001374878035816041),
    @("ECs",    "Thbs1", "Tnfrsf11b", "FAPs", 3, 1, 15.18897966666667,  45.566939,   0.01327029680642083, 0.01327029680642082, 3, 1,                  6.93109,              20.79327, 0.9472576624920083,  0.9472576624920083,  105.2761850778367,  947.4856657005299,  0.01257039033342535,  0.01257039033342535),
    @("ECs",    "Thbs1", "Tnfrsf11b", "sCs",  3, 1, 15.18897966666667,  45.566939,   0.01327029680642083, 0.01327029680642082, 3, 1,                  0.3101076666666667,  0.930323, 0.04238177017576132, 0.04238177017576132, 4.710219043477444,  42.391971391297,    0.0005624186694138669,0.0005624186694138666),

    @("FAPs",   "Thbs1", "Tnfrsf11b", "ECs",  3, 1, 95.50314333333334,  286.50943,   0.08343911742543099, 0.08343911742543098, 1, 0.3333333333333333, 0.07580833333333332, 0.227425, 0.01036056733223033, 0.01036056733223033, 7.239934124194444,  65.15940711774999,  0.0008644765942280508,0.0008644765942280508),
    @("FAPs",   "Thbs1", "Tnfrsf11b", "FAPs", 3, 1, 95.50314333333334,  286.50943,   0.08343911742543099, 0.08343911742543098, 3, 1,                  6.93109,              20.79327, 0.9472576624920083,  0.9472576624920083,  661.9408817262334,  5957.4679355361,    0.07903834333280996,  0.07903834333280994),
    @("FAPs",   "Thbs1", "Tnfrsf11b", "sCs",  3, 1, 95.50314333333334,  286.50943,   0.08343911742543099, 0.08343911742543098, 3, 1,                  0.3101076666666667,  0.930323, 0.04238177017576132, 0.04238177017576132, 29.61625693843223,  266.54631244589,    0.003536297498392978, 0.003536297498392977),

    @("M1",     "Thbs1", "Tnfrsf11b", "ECs",  3, 1, 402.0913113333333,  1206.273934, 0.3512988470441011,  0.351298847044101,   1, 0.3333333333333333, 0.07580833333333332, 0.227425, 0.01036056733223033, 0.01036056733223033, 30.48187215999444,  274.33684943995,    0.003639655358535293, 0.003639655358535292),
    @("M1",     "Thbs1", "Tnfrsf11b", "FAPs", 3, 1, 402.0913113333333,  1206.273934, 0.3512988470441011,  0.351298847044101,   3, 1,                  6.93109,              20.79327, 0.9472576624920083,  0.9472576624920083,  2786.931067069353,  25082.37960362418,  0.3327705246871328,   0.3327705246871326),
    @("M1",     "Thbs1", "Tnfrsf11b", "sCs",  3, 1, 402.0913113333333,  1206.273934, 0.3512988470441011,  0.351298847044101,   3, 1,                  0.3101076666666667,  0.930323, 0.04238177017576132, 0.04238177017576132, 124.6915983445202,  1122.224385100682,  0.01488866699843302,  0.01488866699843302),

    @("M2",     "Thbs1", "Tnfrsf11b", "ECs",  3, 1, 192.2123056666667,  576.636917,  0.1679319086614409,  0.1679319086614409,  1, 0.3333333333333333, 0.07580833333333332, 0.227425, 0.01036056733223033, 0.01036056733223033, 14.57129453874722,  131.141650848725,   0.001739869846916812, 0.001739869846916813),
    @("M2",     "Thbs1", "Tnfrsf11b", "FAPs", 3, 1, 192.2123056666667,  576.636917,  0.1679319086614409,  0.1679319086614409,  3, 1,                  6.93109,              20.79327, 0.9472576624920083,  0.9472576624920083,  1332.240789683177,  11990.16710714859,  0.159074787256458,    0.159074787256458),
    @("M2",     "Thbs1", "Tnfrsf11b", "sCs",  3, 1, 192.2123056666667,  576.636917,  0.1679319086614409,  0.1679319086614409,  3, 1,                  0.3101076666666667,  0.930323, 0.04238177017576132, 0.04238177017576132, 59.60650961491012,  536.458586534191,   0.007117251558066131, 0.007117251558066131),

    @("Neutro", "Thbs1", "Tnfrsf11b", "ECs",  3, 1, 402.444082,         1207.332246, 0.3516070554658648,  0.3516070554658648,  1, 0.3333333333333333, 0.07580833333333332, 0.227425, 0.01036056733223033, 0.01036056733223033, 30.50861511628333,  274.57753604655,    0.003642848572641337, 0.003642848572641337),
    @("Neutro", "Thbs1", "Tnfrsf11b", "FAPs", 3, 1, 402.444082,         1207.332246, 0.3516070554658648,  0.3516070554658648,  3, 1,                  6.93109,              20.79327, 0.9472576624920083,  0.9472576624920083,  2789.37615230938,   25104.38537078442,  0.333062477476293,    0.333062477476293),
    @("Neutro", "Thbs1", "Tnfrsf11b", "sCs",  3, 1, 402.444082,         1207.332246, 0.3516070554658648,  0.3516070554658648,  3, 1,                  0.3101076666666667,  0.930323, 0.04238177017576132, 0.04238177017576132, 124.8009952328287,  1123.208957095458,  0.01490172941693045,  0.01490172941693045),

    @("sCs",    "Thbs1", "Tnfrsf11b", "ECs",  3, 1, 37.14495166666666,  111.434855,  0.03245277459674146, 0.03245277459674145, 1, 0.3333333333333333, 0.07580833333333332, 0.227425, 0.01036056733223033, 0.01036056733223033, 2.815896877597222,  25.343071898375,    0.0003362291563272339,0.0003362291563272338),
    @("sCs",    "Thbs1", "Tnfrsf11b", "FAPs", 3, 1, 37.14495166666666,  111.434855,  0.03245277459674146, 0.03245277459674145, 3, 1,                  6.93109,              20.79327, 0.9472576624920083,  0.9472576624920083,  257.4550030473167,  2317.09502742585,   0.03074113940588934,  0.03074113940588933),
    @("sCs",    "Thbs1", "Tnfrsf11b", "sCs",  3, 1, 37.14495166666666,  111.434855,  0.03245277459674146, 0.03245277459674145, 3, 1,                  0.3101076666666667,  0.930323, 0.04238177017576132, 0.04238177017576132, 11.51893428979611,  103.670408608165,   0.001375406034524882, 0.001375406034524881)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
